$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# ORGANIZATION sheet - add 5 new hospital rows (rows 4-8) beneath the
# existing MOH / Dau Tieng rows, matching the ACC_CREATE_ORG data load.
# Cell-by-cell order below is chosen so that the generated shared-string
# table and cell style table come out in the same order/shape as the
# original author's save.
# ===========================================================================

# ---- Column A : ID (text, formatted as Text "@" so leading zeros show) ----
$ws.Range("A4").Value = "'02001"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "02002"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "'02003"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "02004"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "'02005"
# Existing numeric ID in A3 also gets the Text number format applied.
$ws.Range("A3").NumberFormat = "@"

# ---- Column B : NAME ----
$ws.Range("B4").Value = "BV Đa khoa quốc tế Vinmec Central Park"
$ws.Range("B5").Value = "BV Chợ Rẫy"
$ws.Range("B6").Value = "BV Nhân Dân 115"
$ws.Range("B7").Value = "BV Y Dược TPHCM"
$ws.Range("B8").Value = "BV Nhi Đồng 1"

# ---- Column C : PROVINCENAME ----
$ws.Range("C4").Value = "Hồ Chí Minh"
$ws.Range("C5").Value = "Hồ Chí Minh"
$ws.Range("C6").Value = "Hồ Chí Minh"
$ws.Range("C7").Value = "Hồ Chí Minh"
$ws.Range("C8").Value = "Hồ Chí Minh"

# ---- Column D : DISTRICTNAME ----
$ws.Range("D4").Value = "Bình Thạnh"
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 10

# ---- Column E : TOWNNAME ----
$ws.Range("E4").Value = 22
$ws.Range("E5").Value = 12
$ws.Range("E6").Value = 12
$ws.Range("E7").Value = 11
$ws.Range("E8").Value = 10

# ---- Column F : STREET ----
$ws.Range("F4").Value = "208 Nguyễn Hữu Cảnh"
$ws.Range("F6").Value = "527 Sư Vạn Hạnh"
$ws.Range("F7").Value = "215 Hồng Bàng"
$ws.Range("F8").Value = "341 Sư Vạn Hạnh"
$ws.Range("F5").Value = "201B Nguyễn Chí Thanh"

# ---- Column widths / print setup ----
$ws.Columns("B").ColumnWidth = 36.1666
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
